$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$nl = [char]10

$ws.Range("B10").Value = "< 10 Manuf., Agriculture, and Forestry" + $nl + "< 5 Wholesale, and Retail" + $nl + "< 10 Services, and Mining"
$ws.Range("C10").Value = "< CFA 20 Millionlion Manuf., Agriculture, and Forestry" + $nl + "< CFA 15 Millionlion Wholesale" + $nl + "< CFA 10 Millionlion Retail, Services, and Mining"

$ws.Range("B11").Value = "< 50 Manuf., Agriculture, Forestry, and Wholesale" + $nl + "< 30 Retail, Services, and Mining"
$ws.Range("C11").Value = "< CFA 250 Millionlion Manuf., Agriculture, and Forestry " + $nl + "< CFA 150 Millionlion Wholesale" + $nl + "< CFA 50 Millionlion Retail" + $nl + "< CFA 75 Millionlion Services, and Mining"
$ws.Range("D11").Value = "< CFA 250 Millionlion Manuf., Agriculture, and Forestry " + $nl + "< CFA 200 Millionlion Wholesale" + $nl + "< CFA 100 Millionlion Retail, Services, and Mining"

$ws.Range("B12").Value = "< 100 Manuf., Agriculture, and Forestry" + $nl + "< 50 Wholesale, Retail, Services, and Mining"
$ws.Range("C12").Value = "< CFA 750 Millionlion Manuf., Agriculture, and Forestry " + $nl + "< CFA 250 Millionlion Wholesale, and Retail, Services, and Mining"
$ws.Range("D12").Value = "< CFA 500 Millionlion Manuf., Agriculture, and Forestry " + $nl + "< CFA 250 Millionlion Wholesale, and Retail, Services, and Mining"

$ws.Range("B13").Value = "> 100 Manuf., Agriculture, and Forestry" + $nl + "> 50 Wholesale, Retail, Services, and Mining"
$ws.Range("C13").Value = "> CFA 750 Millionlion Manuf., Agriculture, and Forestry " + $nl + "> CFA 250 Millionlion Wholesale, and Retail, Services, and Mining"
$ws.Range("D13").Value = "> CFA 500 Millionlion Manuf., Agriculture, and Forestry " + $nl + "> CFA 250 Millionlion Wholesale, and Retail, Services, and Mining"
